$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# A new "Late" related column is being introduced for Variable Instalments / RBI
# loan schedules. Insert a blank column before the existing "Late" column (N),
# pushing "Late" to O and "Outstanding" to Q, leaving the new N column blank.
$ws.Columns.Item(14).Insert()

# Match the column width Excel applied to the freshly inserted column.
$ws.Columns.Item(14).ColumnWidth = 9.166666666666666

# Restore sheet activation / selected cell state as recorded in the workbook.
$ws.Activate() | Out-Null
$ws.Range("T6").Select() | Out-Null
